# Update cryptos list values per scraped diff (Sat Aug 12 06:41:41 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.375.86"
$ws.Range("D3").Value = "1.847.09"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6298"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07591"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07744"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "1.843.92"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.001"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001084"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.06%  "
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").Value = "2.088.09"
$ws.Range("E17").Value = "  -7.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.150"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").Value = "29.404.83"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.422"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.06%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1390"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.384"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.313"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.464"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05602"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.029"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.846"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.155"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7091"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.583"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("D38").Value = "1.233.99"
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.771"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01799"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.449"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9075"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("E46").Value = "  +3.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.192"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4014"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.680"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.95%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.955"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.30%  "
$ws.Range("E51").Value = "  -0.72%  "